# Add an "Item" / "Quantity" table (a Macbook pro order line) next to the
# existing TestCaseName / Execution columns on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: "Item" header + "Macbook pro" value
$ws.Range("C1").Value = "Item"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C2").Value = "'Macbook pro"

# Column D: "Quantity" header + "2" value (entered with a leading apostrophe,
# so it is stored as text rather than a number)
$ws.Range("D1").Value = "Quantity"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D2").Value = "'2"

# Resize column C to fit its new (wider) contents
$ws.Range("C1").ColumnWidth = 11

# Leave the selection where the editor ended up
[void]$ws.Range("M15").Select()
